# Regenerate save_data: replace column G ("K" = strikeouts, formerly "Strike#")
# values with the correct strikeout counts for each outing (rows 2-51).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 2
    3  = 3
    4  = 1
    5  = 1
    6  = 0
    7  = 2
    8  = 0
    9  = 1
    10 = 0
    11 = 1
    12 = 0
    13 = 1
    14 = 0
    15 = 0
    16 = 0
    17 = 0
    18 = 1
    19 = 0
    20 = 1
    21 = 0
    22 = 1
    23 = 1
    24 = 1
    25 = 2
    26 = 1
    27 = 2
    28 = 1
    29 = 0
    30 = 0
    31 = 1
    32 = 0
    33 = 2
    34 = 1
    35 = 1
    36 = 2
    37 = 2
    38 = 2
    39 = 1
    40 = 1
    41 = 1
    42 = 1
    43 = 1
    44 = 1
    45 = 2
    46 = 1
    47 = 2
    48 = 3
    49 = 1
    50 = 0
    51 = 0
}

foreach ($row in $newK.Keys) {
    $ws.Cells.Item($row, 7).Value = $newK[$row]
}
